$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column B (old B/C shift to F/G)
$ws.Range("B1:E1").EntireColumn.Insert() | Out-Null

# Row 1 - headers
$ws.Range("B1").Value = "apellido1"
$ws.Range("C1").Value = "apellido2"
$ws.Range("D1").Value = "brigada"
$ws.Range("E1").Value = "especialidad"

# Rows 2-3 for columns B (apellido1) and C (apellido2)
$ws.Range("B2").Value = "martinez"
$ws.Range("C2").Value = "soriano"
$ws.Range("B3").Value = "hernandez"
$ws.Range("C3").Value = "de zuloaga"

# Column D (brigada) - numeric
$ws.Range("D2").Value = 400
$ws.Range("D3").Value = 401

# Column E (especialidad)
$ws.Range("E2").Value = "oss"
$ws.Range("E3").Value = "ham"

# Column widths for new columns (match target: width ~15.71, customWidth, no bestFit)
$ws.Range("B1:E1").EntireColumn.ColumnWidth = 14.83

# Update selection
$ws.Range("E13").Select() | Out-Null
